# Updates cryptos list (GitHub Actions style refresh of Price / Volume(1h)
# columns, plus one coin swap in row 51). Cells in column D whose new value
# looks like a plain number (e.g. "702.78") are explicitly formatted as text
# ("@") before the value is written so Excel keeps them as text strings
# (matching the workbook's existing inline-string cells) instead of silently
# converting them to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.364.80"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "3.776.93"
$ws.Range("E3").Value = "  -2.65%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "702.78"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.39"
$ws.Range("E6").Value = "  -2.69%  "
$ws.Range("D7").Value = "3.779.96"
$ws.Range("E7").Value = "  -2.56%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.518"
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("E12").Value = "  -2.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("E13").Value = "  -3.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.88"
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("D15").Value = "4.410.36"
$ws.Range("E15").Value = "  -2.64%  "
$ws.Range("D16").Value = "3.722.01"
$ws.Range("E16").Value = "  -3.89%  "
$ws.Range("D17").Value = "70.339.39"
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.09"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.18"
$ws.Range("E20").Value = "  -3.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "489.39"
$ws.Range("E21").Value = "  -1.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.51"
$ws.Range("E22").Value = "  -5.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.720"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.69"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("E25").Value = "  -3.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.99"
$ws.Range("E26").Value = "  -3.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.37"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("D28").Value = "3.926.02"
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.03"
$ws.Range("E30").Value = "  -5.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.06"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.27"
$ws.Range("E32").Value = "  -5.16%  "
$ws.Range("E33").Value = "  -4.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.95"
$ws.Range("E34").Value = "  -3.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.173"
$ws.Range("E35").Value = "  -2.83%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "3.743.67"
$ws.Range("E37").Value = "  -2.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.97"
$ws.Range("E38").Value = "  -4.16%  "
$ws.Range("E39").Value = "  -3.60%  "
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.28"
$ws.Range("E41").Value = "  -4.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.88"
$ws.Range("E42").Value = "  -3.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.24"
$ws.Range("E43").Value = "  -5.90%  "
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.89"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.73"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "415.37"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.63"
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.36"
$ws.Range("E51").Value = "  -2.62%  "
